$d = $word.ActiveDocument

# The original run in the 4th paragraph reads:
#   "9号开始第一天去部门报到，"
# The target splits it into two runs (identical formatting):
#   "9号开始第一天去部门"
#   "<expanded text>"
# First perform a plain text replacement so the paragraph ends up with the
# full, expanded wording in a single run.
$newTail = "就差点迟到，因为不认识路，找工位费了点时间。当时还是在E2-5F-1岛，第一次晨会的时候听到大家讨论的内容感觉既兴奋又紧张。兴奋是因为终于加入了紧贴科技最前沿的手机行业，作为数码迷的我喜不自胜。紧张是因为本次跳槽跨行业，对驱动知之甚少，担心跟不上节奏。"

$d.Content.Find.Execute("9号开始第一天去部门报到，", $true, $false, $false, $false, $false, $true, 1, $false, "9号开始第一天去部门" + $newTail, 2)

# Now split that single run into two separate runs (with identical rPr) right
# after "部门", matching the diff's run layout. Splitting the paragraph in two
# and immediately rejoining it (by deleting the intervening paragraph mark)
# leaves the text on either side of the split point in distinct <w:r>
# elements instead of getting re-coalesced into one run.
$full = $d.Content
$splitPos = $full.Text.IndexOf("9号开始第一天去部门") + 10

$insertionPoint = $d.Range($splitPos, $splitPos)
$insertionPoint.InsertParagraphAfter()

$paraMark = $d.Range($splitPos, $splitPos + 1)
$paraMark.Delete()
